$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '43.783.70'
$ws.Cells.Item(2, 5).Value = '  -0.14%  '
$ws.Cells.Item(3, 4).Value = '2.344.24'
$ws.Cells.Item(3, 5).Value = '  -0.49%  '
$ws.Cells.Item(4, 5).Value = '  -0.30%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = '@'
$cell.Value = '239.15'
$cell.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -1.34%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.667'
$cell.Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -3.65%  '
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = '@'
$cell.Value = '72.41'
$cell.Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  -5.97%  '
$ws.Cells.Item(8, 5).Value = '  -0.08%  '
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.602'
$cell.Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  -5.03%  '
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0995'
$cell.Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -2.87%  '
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = '@'
$cell.Value = '57.97'
$cell.Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +1.06%  '
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = '@'
$cell.Value = '32.66'
$cell.Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  -3.17%  '
$ws.Cells.Item(13, 5).Value = '  -0.71%  '
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = '@'
$cell.Value = '7.23'
$cell.Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -4.93%  '
$ws.Cells.Item(15, 4).Value = '2.688.47'
$ws.Cells.Item(15, 5).Value = '  -0.87%  '
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = '@'
$cell.Value = '16.09'
$cell.Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  -5.62%  '
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.902'
$cell.Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  -2.80%  '
$ws.Cells.Item(18, 4).Value = '2.340.95'
$ws.Cells.Item(18, 5).Value = '  -0.53%  '
$ws.Cells.Item(19, 4).Value = '43.734.32'
$ws.Cells.Item(19, 5).Value = '  -0.05%  '
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0000101'
$cell.Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -2.06%  '
$ws.Cells.Item(21, 2).Value = 'Uniswap'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.62'
$cell.Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  -1.10%  '
$ws.Cells.Item(22, 2).Value = 'Litecoin'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = '@'
$cell.Value = '77.74'
$cell.Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -0.01%  '
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = '@'
$cell.Value = '251.73'
$cell.Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -2.08%  '
$ws.Cells.Item(24, 5).Value = '  +0.09%  '
$ws.Cells.Item(25, 2).Value = 'WEMIXToken'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = '@'
$cell.Value = '3.71'
$cell.Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +2.22%  '
$ws.Cells.Item(26, 2).Value = 'ImmutableX'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.85'
$cell.Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  +3.29%  '
$ws.Cells.Item(27, 5).Value = '  -2.05%  '
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = '@'
$cell.Value = '10.33'
$cell.Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  -6.49%  '
$ws.Cells.Item(29, 5).Value = '  -1.29%  '
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = '@'
$cell.Value = '176.31'
$cell.Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  +0.57%  '
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = '@'
$cell.Value = '22.21'
$cell.Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -3.62%  '
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.125'
$cell.Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -2.79%  '
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.133'
$cell.Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -1.97%  '
$ws.Cells.Item(34, 5).Value = '  -3.55%  '
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = '@'
$cell.Value = '5.09'
$cell.Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  -5.14%  '
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = '@'
$cell.Value = '5.33'
$cell.Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  -1.14%  '
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = '@'
$cell.Value = '3.74'
$cell.Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -1.14%  '
$ws.Cells.Item(38, 2).Value = 'FTXToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = '@'
$cell.Value = '5.86'
$cell.Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  +30.41%  '
$ws.Cells.Item(39, 2).Value = 'THORChain'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.37'
$cell.Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -1.56%  '
$ws.Cells.Item(40, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.37'
$cell.Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -2.81%  '
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0272'
$cell.Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -2.58%  '
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = '@'
$cell.Value = '66.17'
$cell.Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  +15.59%  '
$ws.Cells.Item(43, 2).Value = 'FraxShare'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = '@'
$cell.Value = '9.18'
$cell.Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +1.22%  '
$ws.Cells.Item(44, 2).Value = 'Cronos'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.107'
$cell.Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  +1.62%  '
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = '@'
$cell.Value = '18.80'
$cell.Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -3.12%  '
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.196'
$cell.Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -3.72%  '
$ws.Cells.Item(47, 5).Value = '  -0.14%  '
$ws.Cells.Item(48, 5).Value = '  -3.64%  '
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.42'
$cell.Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -4.79%  '
$ws.Cells.Item(50, 2).Value = 'ARBITRUM'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.14'
$cell.Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -4.14%  '
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = '@'
$cell.Value = '97.97'
$cell.Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  -3.94%  '
